$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data occupies A3:K12. Pre-format as Text ("@") so the numeric-looking
# strings (e.g. "106", "212.00") are stored as text, matching the existing
# dataset (columns G:K keep trailing zeros / leading content exactly).
$dataRange = $ws.Range("A3:K12")
$dataRange.NumberFormat = "@"

# Row 3
$ws.Cells.Item(3, 1).Value = ' Sep 27 2020'
$ws.Cells.Item(3, 2).Value = ' Sharjah'
$ws.Cells.Item(3, 3).Value = 'Royals won by 4 wickets (with 3 balls remaining)'
$ws.Cells.Item(3, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(3, 5).Value = 'Rajasthan Royals'
$ws.Cells.Item(3, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(3, 7).Value = '106'
$ws.Cells.Item(3, 8).Value = '50'
$ws.Cells.Item(3, 9).Value = '10'
$ws.Cells.Item(3, 10).Value = '7'
$ws.Cells.Item(3, 11).Value = '212.00'

# Row 4
$ws.Cells.Item(4, 1).Value = ' Oct 15 2020'
$ws.Cells.Item(4, 2).Value = ' Sharjah'
$ws.Cells.Item(4, 3).Value = 'Kings XI won by 8 wickets'
$ws.Cells.Item(4, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(4, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(4, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(4, 7).Value = '45'
$ws.Cells.Item(4, 8).Value = '25'
$ws.Cells.Item(4, 9).Value = '4'
$ws.Cells.Item(4, 10).Value = '3'
$ws.Cells.Item(4, 11).Value = '180.00'

# Row 5
$ws.Cells.Item(5, 1).Value = ' Oct 20 2020'
$ws.Cells.Item(5, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(5, 3).Value = 'Kings XI won by 5 wickets (with 6 balls remaining)'
$ws.Cells.Item(5, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(5, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(5, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(5, 7).Value = '5'
$ws.Cells.Item(5, 8).Value = '9'
$ws.Cells.Item(5, 9).Value = '0'
$ws.Cells.Item(5, 10).Value = '0'
$ws.Cells.Item(5, 11).Value = '55.55'

# Row 6
$ws.Cells.Item(6, 1).Value = ' Sep 24 2020'
$ws.Cells.Item(6, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(6, 3).Value = 'Kings XI won by 97 runs'
$ws.Cells.Item(6, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(6, 5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(6, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(6, 7).Value = '26'
$ws.Cells.Item(6, 8).Value = '20'
$ws.Cells.Item(6, 9).Value = '4'
$ws.Cells.Item(6, 10).Value = '0'
$ws.Cells.Item(6, 11).Value = '130.00'

# Row 7
$ws.Cells.Item(7, 1).Value = ' Oct 18 2020'
$ws.Cells.Item(7, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(7, 3).Value = 'Match tied (Kings XI won the one-over eliminator)'
$ws.Cells.Item(7, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(7, 5).Value = 'Mumbai Indians'
$ws.Cells.Item(7, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(7, 7).Value = '11'
$ws.Cells.Item(7, 8).Value = '10'
$ws.Cells.Item(7, 9).Value = '1'
$ws.Cells.Item(7, 10).Value = '0'
$ws.Cells.Item(7, 11).Value = '110.00'

# Row 8
$ws.Cells.Item(8, 1).Value = ' Oct 4 2020'
$ws.Cells.Item(8, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(8, 3).Value = 'Super Kings won by 10 wickets (with 14 balls remaining)'
$ws.Cells.Item(8, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(8, 5).Value = 'Chennai Super Kings'
$ws.Cells.Item(8, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(8, 7).Value = '26'
$ws.Cells.Item(8, 8).Value = '19'
$ws.Cells.Item(8, 9).Value = '3'
$ws.Cells.Item(8, 10).Value = '0'
$ws.Cells.Item(8, 11).Value = '136.84'

# Row 9
$ws.Cells.Item(9, 1).Value = ' Oct 10 2020'
$ws.Cells.Item(9, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(9, 3).Value = 'KKR won by 2 runs'
$ws.Cells.Item(9, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(9, 5).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(9, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(9, 7).Value = '56'
$ws.Cells.Item(9, 8).Value = '39'
$ws.Cells.Item(9, 9).Value = '6'
$ws.Cells.Item(9, 10).Value = '1'
$ws.Cells.Item(9, 11).Value = '143.58'

# Row 10
$ws.Cells.Item(10, 1).Value = ' Oct 8 2020'
$ws.Cells.Item(10, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(10, 3).Value = 'Sunrisers won by 69 runs'
$ws.Cells.Item(10, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(10, 5).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(10, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(10, 7).Value = '9'
$ws.Cells.Item(10, 8).Value = '6'
$ws.Cells.Item(10, 9).Value = '1'
$ws.Cells.Item(10, 10).Value = '0'
$ws.Cells.Item(10, 11).Value = '150.00'

# Row 11
$ws.Cells.Item(11, 1).Value = ' Sep 20 2020'
$ws.Cells.Item(11, 2).Value = ' Dubai (DSC)'
$ws.Cells.Item(11, 3).Value = 'Match tied (Capitals won the one-over eliminator)'
$ws.Cells.Item(11, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(11, 5).Value = 'Delhi Capitals'
$ws.Cells.Item(11, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(11, 7).Value = '89'
$ws.Cells.Item(11, 8).Value = '60'
$ws.Cells.Item(11, 9).Value = '7'
$ws.Cells.Item(11, 10).Value = '4'
$ws.Cells.Item(11, 11).Value = '148.33'

# Row 12
$ws.Cells.Item(12, 1).Value = ' Oct 1 2020'
$ws.Cells.Item(12, 2).Value = ' Abu Dhabi'
$ws.Cells.Item(12, 3).Value = 'Mumbai won by 48 runs'
$ws.Cells.Item(12, 4).Value = 'Kings XI Punjab'
$ws.Cells.Item(12, 5).Value = 'Mumbai Indians'
$ws.Cells.Item(12, 6).Value = 'Mayank Agarwal '
$ws.Cells.Item(12, 7).Value = '25'
$ws.Cells.Item(12, 8).Value = '18'
$ws.Cells.Item(12, 9).Value = '3'
$ws.Cells.Item(12, 10).Value = '0'
$ws.Cells.Item(12, 11).Value = '138.88'

# Drop the explicit Text number-format style now that the values are locked
# in as text, so the new cells end up on the default style like the rest of
# the sheet.
$dataRange.Style = "Normal"
